$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") values are plain text in the source workbook (e.g. "1.00",
# "27.912.30" with a thousands separator). Excel auto-converts Value assignments
# that look numeric into real numbers, so force the cell to Text format first,
# then reset the style back to Normal so no stray number-format style lingers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.912.30"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.01%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.636.88"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("E6").Value = "  -0.57%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.42"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.77%  "
$ws.Range("E9").Value = "  -0.33%  "
$ws.Range("E10").Value = "  -0.28%  "
$ws.Range("E11").Value = "  +0.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.868.20"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.632.38"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.00%  "
$ws.Range("E14").Value = "  -0.95%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.564"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.37"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.15%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.925.57"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "229.48"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.81"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.03%  "
$ws.Range("E20").Value = "  -0.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.999"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.83%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.65%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "156.28"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.28%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("E28").Value = "  -0.36%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0482"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.42"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.12"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.403.89"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.36%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.61"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.07%  "
$ws.Range("E36").Value = "  +0.96%  "
$ws.Range("E37").Value = "  -0.75%  "
$ws.Range("E38").Value = "  +0.54%  "
$ws.Range("E39").Value = "  -0.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.854"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E42").Value = "  -1.06%  "
$ws.Range("E43").Value = "  +2.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "66.21"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.89%  "
$ws.Range("E45").Value = "  -1.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.776.76"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.72"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.22%  "
$ws.Range("E49").Value = "  +2.06%  "
$ws.Range("E50").Value = "  -0.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.62"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.81%  "
